# Mark attendance-related cells with value 1 for specific (row, column) pairs
# as recorded by the attendance tracking macro for this student's sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("G3", "H3", "H4", "H5", "D6", "E6", "H7", "D8", "E8", "H9", "H10", "H11", "D12", "E12", "H13", "H14", "H15", "H16", "H17", "H18")

foreach ($cellRef in $cells) {
    $ws.Range($cellRef).Value = 1
}
